$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 17.52333086377223
$ws.Cells.Item(2, 3).Value = 10.01068817667626
$ws.Cells.Item(2, 4).Value = 6.038340297485728
$ws.Cells.Item(2, 5).Value = 16.1700202647893
$ws.Cells.Item(2, 7).Value = 3.719602636751555
$ws.Cells.Item(2, 11).Value = 16.87841560734294
$ws.Cells.Item(2, 14).Value = 24.48014861058435

$ws.Cells.Item(3, 2).Value = 17.18910987679533
$ws.Cells.Item(3, 3).Value = 9.711022583902013
$ws.Cells.Item(3, 4).Value = 5.930281471181327
$ws.Cells.Item(3, 5).Value = 15.27057575901039
$ws.Cells.Item(3, 7).Value = 3.724140444865302
$ws.Cells.Item(3, 11).Value = 16.63006494797013
$ws.Cells.Item(3, 14).Value = 24.36685345874635

$ws.Cells.Item(4, 2).Value = 16.98661358983785
$ws.Cells.Item(4, 3).Value = 9.5260670179636
$ws.Cells.Item(4, 4).Value = 5.864968902130087
$ws.Cells.Item(4, 5).Value = 14.69649287162362
$ws.Cells.Item(4, 7).Value = 3.727063145343548
$ws.Cells.Item(4, 11).Value = 16.4812683904895
$ws.Cells.Item(4, 14).Value = 24.2980326955226

$ws.Cells.Item(5, 2).Value = 16.90489598540878
$ws.Cells.Item(5, 3).Value = 9.450584201949361
$ws.Cells.Item(5, 4).Value = 5.838651050039871
$ws.Cells.Item(5, 5).Value = 14.45734178204565
$ws.Cells.Item(5, 7).Value = 3.728288647617414
$ws.Cells.Item(5, 11).Value = 16.42163519218376
$ws.Cells.Item(5, 14).Value = 24.27018153891894

$ws.Cells.Item(6, 2).Value = 16.89137867773122
$ws.Cells.Item(6, 3).Value = 9.43804748682853
$ws.Cells.Item(6, 4).Value = 5.8343000376142
$ws.Cells.Item(6, 5).Value = 14.4173255914185
$ws.Cells.Item(6, 7).Value = 3.728494228472524
$ws.Cells.Item(6, 11).Value = 16.411795827087
$ws.Cells.Item(6, 14).Value = 24.26556878293903

$ws.Cells.Item(7, 2).Value = 16.98550811749289
$ws.Cells.Item(7, 3).Value = 9.525049307655282
$ws.Cells.Item(7, 4).Value = 5.864612716901958
$ws.Cells.Item(7, 5).Value = 14.69328827462067
$ws.Cells.Item(7, 7).Value = 3.727079533079202
$ws.Cells.Item(7, 11).Value = 16.48046000025251
$ws.Cells.Item(7, 14).Value = 24.29765629110358

$ws.Cells.Item(8, 2).Value = 17.40759850579231
$ws.Cells.Item(8, 3).Value = 9.907638184388153
$ws.Cells.Item(8, 4).Value = 6.000888895104457
$ws.Cells.Item(8, 5).Value = 15.86457257569413
$ws.Cells.Item(8, 7).Value = 3.721139047028846
$ws.Cells.Item(8, 11).Value = 16.79206326766036
$ws.Cells.Item(8, 14).Value = 24.44093457156859

$ws.Cells.Item(9, 2).Value = 18.25170223468708
$ws.Cells.Item(9, 3).Value = 10.64496013628211
$ws.Cells.Item(9, 4).Value = 6.274713732922391
$ws.Cells.Item(9, 5).Value = 17.97847481887194
$ws.Cells.Item(9, 7).Value = 3.710565163346881
$ws.Cells.Item(9, 11).Value = 17.42905880918163
$ws.Cells.Item(9, 14).Value = 24.72747950322624

$ws.Cells.Item(10, 2).Value = 18.87505256207881
$ws.Cells.Item(10, 3).Value = 11.17223697823084
$ws.Cells.Item(10, 4).Value = 6.477745507917873
$ws.Cells.Item(10, 5).Value = 19.56736473808277
$ws.Cells.Item(10, 7).Value = 3.703441692710354
$ws.Cells.Item(10, 11).Value = 17.90832696514771
$ws.Cells.Item(10, 14).Value = 24.94096710015651

$ws.Cells.Item(11, 2).Value = 19.15794360947391
$ws.Cells.Item(11, 3).Value = 11.40776394696344
$ws.Cells.Item(11, 4).Value = 6.570068962762825
$ws.Cells.Item(11, 5).Value = 20.26168117415078
$ws.Cells.Item(11, 7).Value = 3.700338885139817
$ws.Cells.Item(11, 11).Value = 18.12784155151594
$ws.Cells.Item(11, 14).Value = 25.03863261766095

$ws.Cells.Item(12, 2).Value = 19.26486111896585
$ws.Cells.Item(12, 3).Value = 11.49623984815805
$ws.Cells.Item(12, 4).Value = 6.604988830364392
$ws.Cells.Item(12, 5).Value = 20.51884627912217
$ws.Cells.Item(12, 7).Value = 3.699183557565267
$ws.Cells.Item(12, 11).Value = 18.21110170934368
$ws.Cells.Item(12, 14).Value = 25.07568586646896

$ws.Cells.Item(13, 2).Value = 19.24184543326759
$ws.Cells.Item(13, 3).Value = 11.47721802991182
$ws.Cells.Item(13, 4).Value = 6.597470589229462
$ws.Cells.Item(13, 5).Value = 20.46371663928199
$ws.Cells.Item(13, 7).Value = 3.69943150719468
$ws.Cells.Item(13, 11).Value = 18.19316540351722
$ws.Cells.Item(13, 14).Value = 25.06770285190089

$ws.Cells.Item(14, 2).Value = 19.1667445421398
$ws.Cells.Item(14, 3).Value = 11.41505768649067
$ws.Cells.Item(14, 4).Value = 6.572942862671606
$ws.Cells.Item(14, 5).Value = 20.28295347499887
$ws.Cells.Item(14, 7).Value = 3.700243442994337
$ws.Cells.Item(14, 11).Value = 18.13468916520342
$ws.Cells.Item(14, 14).Value = 25.0416796744044

$ws.Cells.Item(15, 2).Value = 19.1207129147453
$ws.Cells.Item(15, 3).Value = 11.37688730521412
$ws.Cells.Item(15, 4).Value = 6.557912545096004
$ws.Cells.Item(15, 5).Value = 20.17148214112984
$ws.Cells.Item(15, 7).Value = 3.700743329785801
$ws.Cells.Item(15, 11).Value = 18.0988859893447
$ws.Cells.Item(15, 14).Value = 25.0257484834095

$ws.Cells.Item(16, 2).Value = 18.85654185646479
$ws.Cells.Item(16, 3).Value = 11.15674964604219
$ws.Cells.Item(16, 4).Value = 6.471708124683924
$ws.Cells.Item(16, 5).Value = 19.5211818065793
$ws.Cells.Item(16, 7).Value = 3.703647225413553
$ws.Cells.Item(16, 11).Value = 17.89400433760668
$ws.Cells.Item(16, 14).Value = 24.93459466851918

$ws.Cells.Item(17, 2).Value = 18.69422870792645
$ws.Cells.Item(17, 3).Value = 11.02052557019733
$ws.Cells.Item(17, 4).Value = 6.418789202156813
$ws.Cells.Item(17, 5).Value = 19.11194863378293
$ws.Cells.Item(17, 7).Value = 3.705463823129836
$ws.Cells.Item(17, 11).Value = 17.76864103084031
$ws.Cells.Item(17, 14).Value = 24.87880867244162

$ws.Cells.Item(18, 2).Value = 18.6008133050431
$ws.Cells.Item(18, 3).Value = 10.94176995927127
$ws.Cells.Item(18, 4).Value = 6.388350231630327
$ws.Cells.Item(18, 5).Value = 18.87277307036648
$ws.Cells.Item(18, 7).Value = 3.706521652307986
$ws.Cells.Item(18, 11).Value = 17.69668059419379
$ws.Cells.Item(18, 14).Value = 24.84677429688144

$ws.Cells.Item(19, 2).Value = 18.56917808037196
$ws.Cells.Item(19, 3).Value = 10.91503830788541
$ws.Cells.Item(19, 4).Value = 6.378044991176278
$ws.Cells.Item(19, 5).Value = 18.79113931175575
$ws.Cells.Item(19, 7).Value = 3.706882047666145
$ws.Cells.Item(19, 11).Value = 17.67234343533235
$ws.Cells.Item(19, 14).Value = 24.83593726769224

$ws.Cells.Item(20, 2).Value = 18.71151395233525
$ws.Cells.Item(20, 3).Value = 11.0350692564058
$ws.Cells.Item(20, 4).Value = 6.424422921004474
$ws.Cells.Item(20, 5).Value = 19.15590474697555
$ws.Cells.Item(20, 7).Value = 3.705269102006263
$ws.Cells.Item(20, 11).Value = 17.78197172424095
$ws.Cells.Item(20, 14).Value = 24.88474184691732

$ws.Cells.Item(21, 2).Value = 19.18880996942706
$ws.Cells.Item(21, 3).Value = 11.43333571046756
$ws.Cells.Item(21, 4).Value = 6.580148641059798
$ws.Cells.Item(21, 5).Value = 20.33620395464268
$ws.Cells.Item(21, 7).Value = 3.700004426083852
$ws.Cells.Item(21, 11).Value = 18.15186202424982
$ws.Cells.Item(21, 14).Value = 25.04932150483961

$ws.Cells.Item(22, 2).Value = 19.49949559685796
$ws.Cells.Item(22, 3).Value = 11.68943413832393
$ws.Cells.Item(22, 4).Value = 6.681669840546749
$ws.Cells.Item(22, 5).Value = 21.07405622120865
$ws.Cells.Item(22, 7).Value = 3.696678047595271
$ws.Cells.Item(22, 11).Value = 18.39435655126976
$ws.Cells.Item(22, 14).Value = 25.15728398651581

$ws.Cells.Item(23, 2).Value = 19.33382666715998
$ws.Cells.Item(23, 3).Value = 11.55316054948319
$ws.Cells.Item(23, 4).Value = 6.627520761204357
$ws.Cells.Item(23, 5).Value = 20.6833072249261
$ws.Cells.Item(23, 7).Value = 3.698442985221364
$ws.Cells.Item(23, 11).Value = 18.26488992793713
$ws.Cells.Item(23, 14).Value = 25.09962896258347

$ws.Cells.Item(24, 2).Value = 18.70369959455941
$ws.Cells.Item(24, 3).Value = 11.02849541502188
$ws.Cells.Item(24, 4).Value = 6.42187596125789
$ws.Cells.Item(24, 5).Value = 19.13604434478799
$ws.Cells.Item(24, 7).Value = 3.7053570935423
$ws.Cells.Item(24, 11).Value = 17.77594456094847
$ws.Cells.Item(24, 14).Value = 24.88205934047805

$ws.Cells.Item(25, 2).Value = 18.02230838325886
$ws.Cells.Item(25, 3).Value = 10.44758903621764
$ws.Cells.Item(25, 4).Value = 6.200157439593871
$ws.Cells.Item(25, 5).Value = 17.42740422750345
$ws.Cells.Item(25, 7).Value = 3.713311627667473
$ws.Cells.Item(25, 11).Value = 17.25442602536067
$ws.Cells.Item(25, 14).Value = 24.64940823883975
